# "Generate Report for Handback"
#
# This script updates the localization-status workbook to reflect that the
# handback files are now in sync with en-US: it records a handback
# timestamp, flips the status text, and adds "Latest Target File" /
# "Latest Handback File" columns (F/G) to the zh-cn and de-de sheets that
# mirror the existing handoff-file hyperlinks.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# ---------------------------------------------------------------------
# 1. Flip the status text everywhere it appears: "Ready for handoff" is
#    now "Handed back: in sync with en-US".
# ---------------------------------------------------------------------
$newStatus = "Handed back: in sync with en-US"

$wsOverview.Range("B2").Value = $newStatus
$wsOverview.Range("C2").Value = $newStatus
$wsOverview.Range("B3").Value = $newStatus
$wsOverview.Range("C3").Value = $newStatus

$wsZhCn.Range("C2").Value = $newStatus
$wsZhCn.Range("C3").Value = $newStatus

$wsDeDe.Range("C2").Value = $newStatus
$wsDeDe.Range("C3").Value = $newStatus

# ---------------------------------------------------------------------
# 2. zh-cn sheet: add "Latest Target File" (F) / "Latest Handback File"
#    (G) columns, reusing the same targets as the existing handoff-file
#    hyperlinks in columns A and D.
# ---------------------------------------------------------------------
$wsZhCn.Hyperlinks.Add(
    $wsZhCn.Range("F2"),
    "https://github.com/OpenLocalizationTest/oltest/blob/e8bc8c0053623ca406bf11dd3cf5ecc039ed5a7a/e2e/aa2e6c52-c2ec-479f-ab23-1d8cd9d808e7.md",
    "",
    "",
    "aa2e6c52-c2ec-479f-ab23-1d8cd9d808e7.md"
) | Out-Null

$wsZhCn.Hyperlinks.Add(
    $wsZhCn.Range("G2"),
    "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/91b9e2bde2cf49e43877f8f36fe9397513e72d1a/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/aa2e6c52-c2ec-479f-ab23-1d8cd9d808e7.1ec808347d77edfe995b1771f407843e40f02d3d.zh-cn.xlf",
    "",
    "",
    "aa2e6c52-c2ec-479f-ab23-1d8cd9d808e7.1ec808347d77edfe995b1771f407843e40f02d3d.zh-cn.xlf"
) | Out-Null

$wsZhCn.Hyperlinks.Add(
    $wsZhCn.Range("F3"),
    "https://github.com/OpenLocalizationTest/oltest/blob/e8bc8c0053623ca406bf11dd3cf5ecc039ed5a7a/e2e/df5ee194-349f-49d0-94b3-ac6dcf090b06.md",
    "",
    "",
    "df5ee194-349f-49d0-94b3-ac6dcf090b06.md"
) | Out-Null

$wsZhCn.Hyperlinks.Add(
    $wsZhCn.Range("G3"),
    "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/91b9e2bde2cf49e43877f8f36fe9397513e72d1a/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/df5ee194-349f-49d0-94b3-ac6dcf090b06.4087c657642ae1b95d210482e79ba24336fbd306.zh-cn.xlf",
    "",
    "",
    "df5ee194-349f-49d0-94b3-ac6dcf090b06.4087c657642ae1b95d210482e79ba24336fbd306.zh-cn.xlf"
) | Out-Null

# Latest Handback DateTime (H) for zh-cn: handback completed.
$wsZhCn.Range("H2").Value = "2016-03-22 11:05:57"
$wsZhCn.Range("H3").Value = "2016-03-22 11:05:57"

# ---------------------------------------------------------------------
# 3. de-de sheet: same idea -- add "Latest Target File" (F) /
#    "Latest Handback File" (G) columns.
# ---------------------------------------------------------------------
$wsDeDe.Hyperlinks.Add(
    $wsDeDe.Range("F2"),
    "https://github.com/OpenLocalizationTest/oltest/blob/e8bc8c0053623ca406bf11dd3cf5ecc039ed5a7a/e2e/aa2e6c52-c2ec-479f-ab23-1d8cd9d808e7.md",
    "",
    "",
    "aa2e6c52-c2ec-479f-ab23-1d8cd9d808e7.md"
) | Out-Null

$wsDeDe.Hyperlinks.Add(
    $wsDeDe.Range("G2"),
    "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/5a00313e8bb789f75b0ecba28a6d8263a0844cdb/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/aa2e6c52-c2ec-479f-ab23-1d8cd9d808e7.1ec808347d77edfe995b1771f407843e40f02d3d.de-de.xlf",
    "",
    "",
    "aa2e6c52-c2ec-479f-ab23-1d8cd9d808e7.1ec808347d77edfe995b1771f407843e40f02d3d.de-de.xlf"
) | Out-Null

$wsDeDe.Hyperlinks.Add(
    $wsDeDe.Range("F3"),
    "https://github.com/OpenLocalizationTest/oltest/blob/e8bc8c0053623ca406bf11dd3cf5ecc039ed5a7a/e2e/df5ee194-349f-49d0-94b3-ac6dcf090b06.md",
    "",
    "",
    "df5ee194-349f-49d0-94b3-ac6dcf090b06.md"
) | Out-Null

$wsDeDe.Hyperlinks.Add(
    $wsDeDe.Range("G3"),
    "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/5a00313e8bb789f75b0ecba28a6d8263a0844cdb/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/df5ee194-349f-49d0-94b3-ac6dcf090b06.4087c657642ae1b95d210482e79ba24336fbd306.de-de.xlf",
    "",
    "",
    "df5ee194-349f-49d0-94b3-ac6dcf090b06.4087c657642ae1b95d210482e79ba24336fbd306.de-de.xlf"
) | Out-Null

# Latest Handback DateTime (H) for de-de: handback completed a little
# later than zh-cn, so it gets its own distinct timestamp.
$wsDeDe.Range("H2").Value = "2016-03-22 11:06:04"
$wsDeDe.Range("H3").Value = "2016-03-22 11:06:04"

$wb.Save()
